$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The underlying "missing data" sample changed: the row for "RM 232" (row 26)
# and the row for "SC 92" (row 28) are no longer part of this extract, so both
# rows are removed and everything below shifts up. Delete the higher-indexed
# row first so the second delete still targets the correct row.
$ws.Rows(28).Delete()
$ws.Rows(26).Delete()

# After the two deletions, row 26 is left holding the old "SC 92" record;
# the new extract instead samples "SC 5" there, so overwrite the whole row.
$ws.Range("A26").Value = "SC 5"
$ws.Range("B26").Value = -20.2
$ws.Range("C26").Value = 10.8
$ws.Range("D26").Value = -13.8
$ws.Range("E26").Value = -5
$ws.Range("F26").Value = 17.38

# After the deletions the sheet holds 32 data rows (A1:F33). Patch the handful
# of cells whose "missing" (blank) vs "present" (value) status differs from a
# plain shift, per the new random missing-data mask.
$ws.Range("F2").Value = 18.03
$ws.Range("F6").Value = ""
$ws.Range("F12").Value = 17.45
$ws.Range("F14").Value = ""
$ws.Range("F20").Value = 17.73
$ws.Range("F21").Value = 16.58
$ws.Range("F23").Value = ""
$ws.Range("F24").Value = ""

$ws.Range("E27").Value = ""
$ws.Range("E28").Value = ""
$ws.Range("E29").Value = -6.8
$ws.Range("E30").Value = -5.7
$ws.Range("E31").Value = ""
$ws.Range("F31").Value = 17.18
$ws.Range("E32").Value = ""
$ws.Range("F33").Value = 17.53
